$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.892066666666667
$ws.Range("H2").Value = 23.6762
$ws.Range("I2").Value = 0.1739002798877711
$ws.Range("J2").Value = 0.1739002798877711
$ws.Range("M2").Value = 1.619868333333333
$ws.Range("N2").Value = 4.859605
$ws.Range("O2").Value = 0.1089327058120143
$ws.Range("P2").Value = 0.1089327058120143
$ws.Range("Q2").Value = 12.78410887788889
$ws.Range("R2").Value = 115.056979901
$ws.Range("S2").Value = 0.01894342802964151
$ws.Range("T2").Value = 0.01894342802964151
$ws.Range("G3").Value = 7.892066666666667
$ws.Range("H3").Value = 23.6762
$ws.Range("I3").Value = 0.1739002798877711
$ws.Range("J3").Value = 0.1739002798877711
$ws.Range("O3").Value = 0.1655705935257241
$ws.Range("P3").Value = 0.1655705935257241
$ws.Range("Q3").Value = 19.4310099876
$ws.Range("R3").Value = 174.8790898884
$ws.Range("S3").Value = 0.0287927725553078
$ws.Range("T3").Value = 0.0287927725553078
$ws.Range("G4").Value = 7.892066666666667
$ws.Range("H4").Value = 23.6762
$ws.Range("I4").Value = 0.1739002798877711
$ws.Range("J4").Value = 0.1739002798877711
$ws.Range("M4").Value = 7.682722666666667
$ws.Range("N4").Value = 23.048168
$ws.Range("O4").Value = 0.5166467859527435
$ws.Range("P4").Value = 0.5166467859527435
$ws.Range("Q4").Value = 60.63255946684445
$ws.Range("R4").Value = 545.6930352016001
$ws.Range("S4").Value = 0.08984502068029944
$ws.Range("T4").Value = 0.08984502068029944
$ws.Range("G5").Value = 7.892066666666667
$ws.Range("H5").Value = 23.6762
$ws.Range("I5").Value = 0.1739002798877711
$ws.Range("J5").Value = 0.1739002798877711
$ws.Range("M5").Value = 3.105673
$ws.Range("N5").Value = 9.317019
$ws.Range("O5").Value = 0.2088499147095181
$ws.Range("P5").Value = 0.2088499147095181
$ws.Range("Q5").Value = 24.51017836086667
$ws.Range("R5").Value = 220.5916052478
$ws.Range("S5").Value = 0.03631905862252231
$ws.Range("T5").Value = 0.03631905862252231
$ws.Range("I6").Value = 0.3815924715300191
$ws.Range("J6").Value = 0.3815924715300191
$ws.Range("M6").Value = 1.619868333333333
$ws.Range("N6").Value = 4.859605
$ws.Range("O6").Value = 0.1089327058120143
$ws.Range("P6").Value = 0.1089327058120143
$ws.Range("Q6").Value = 28.05239707590334
$ws.Range("R6").Value = 252.47157368313
$ws.Range("S6").Value = 0.04156790044125901
$ws.Range("T6").Value = 0.04156790044125901
$ws.Range("I7").Value = 0.3815924715300191
$ws.Range("J7").Value = 0.3815924715300191
$ws.Range("O7").Value = 0.1655705935257241
$ws.Range("P7").Value = 0.1655705935257241
$ws.Range("S7").Value = 0.06318049199617325
$ws.Range("T7").Value = 0.06318049199617325
$ws.Range("I8").Value = 0.3815924715300191
$ws.Range("J8").Value = 0.3815924715300191
$ws.Range("M8").Value = 7.682722666666667
$ws.Range("N8").Value = 23.048168
$ws.Range("O8").Value = 0.5166467859527435
$ws.Range("P8").Value = 0.5166467859527435
$ws.Range("Q8").Value = 133.0471016899787
$ws.Range("R8").Value = 1197.423915209808
$ws.Range("S8").Value = 0.1971485239597481
$ws.Range("T8").Value = 0.1971485239597482
$ws.Range("I9").Value = 0.3815924715300191
$ws.Range("J9").Value = 0.3815924715300191
$ws.Range("M9").Value = 3.105673
$ws.Range("N9").Value = 9.317019
$ws.Range("O9").Value = 0.2088499147095181
$ws.Range("P9").Value = 0.2088499147095181
$ws.Range("Q9").Value = 53.783119523446
$ws.Range("R9").Value = 484.0480757110141
$ws.Range("S9").Value = 0.07969555513283869
$ws.Range("T9").Value = 0.0796955551328387
$ws.Range("G10").Value = 7.716272666666666
$ws.Range("H10").Value = 23.148818
$ws.Range("I10").Value = 0.1700266904854272
$ws.Range("J10").Value = 0.1700266904854272
$ws.Range("M10").Value = 1.619868333333333
$ws.Range("N10").Value = 4.859605
$ws.Range("O10").Value = 0.1089327058120143
$ws.Range("P10").Value = 0.1089327058120143
$ws.Range("Q10").Value = 12.49934574409889
$ws.Range("R10").Value = 112.49411169689
$ws.Range("S10").Value = 0.01852146745483945
$ws.Range("T10").Value = 0.01852146745483945
$ws.Range("G11").Value = 7.716272666666666
$ws.Range("H11").Value = 23.148818
$ws.Range("I11").Value = 0.1700266904854272
$ws.Range("J11").Value = 0.1700266904854272
$ws.Range("O11").Value = 0.1655705935257241
$ws.Range("P11").Value = 0.1655705935257241
$ws.Range("Q11").Value = 18.998188634964
$ws.Range("R11").Value = 170.983697714676
$ws.Range("S11").Value = 0.02815142005888678
$ws.Range("T11").Value = 0.02815142005888678
$ws.Range("G12").Value = 7.716272666666666
$ws.Range("H12").Value = 23.148818
$ws.Range("I12").Value = 0.1700266904854272
$ws.Range("J12").Value = 0.1700266904854272
$ws.Range("M12").Value = 7.682722666666667
$ws.Range("N12").Value = 23.048168
$ws.Range("O12").Value = 0.5166467859527435
$ws.Range("P12").Value = 0.5166467859527435
$ws.Range("Q12").Value = 59.28198291838044
$ws.Range("R12").Value = 533.5378462654239
$ws.Range("S12").Value = 0.08784374316547788
$ws.Range("T12").Value = 0.0878437431654779
$ws.Range("G13").Value = 7.716272666666666
$ws.Range("H13").Value = 23.148818
$ws.Range("I13").Value = 0.1700266904854272
$ws.Range("J13").Value = 0.1700266904854272
$ws.Range("M13").Value = 3.105673
$ws.Range("N13").Value = 9.317019
$ws.Range("O13").Value = 0.2088499147095181
$ws.Range("P13").Value = 0.2088499147095181
$ws.Range("Q13").Value = 23.96421968150466
$ws.Range("R13").Value = 215.677977133542
$ws.Range("S13").Value = 0.0355100598062231
$ws.Range("T13").Value = 0.03551005980622311
$ws.Range("G14").Value = 12.45667266666667
$ws.Range("H14").Value = 37.370018
$ws.Range("I14").Value = 0.2744805580967825
$ws.Range("J14").Value = 0.2744805580967826
$ws.Range("M14").Value = 1.619868333333333
$ws.Range("N14").Value = 4.859605
$ws.Range("O14").Value = 0.1089327058120143
$ws.Range("P14").Value = 0.1089327058120143
$ws.Range("Q14").Value = 20.17816959143223
$ws.Range("R14").Value = 181.60352632289
$ws.Range("S14").Value = 0.02989990988627431
$ws.Range("T14").Value = 0.02989990988627431
$ws.Range("G15").Value = 12.45667266666667
$ws.Range("H15").Value = 37.370018
$ws.Range("I15").Value = 0.2744805580967825
$ws.Range("J15").Value = 0.2744805580967826
$ws.Range("O15").Value = 0.1655705935257241
$ws.Range("P15").Value = 0.1655705935257241
$ws.Range("Q15").Value = 30.669499032564
$ws.Range("R15").Value = 276.025491293076
$ws.Range("S15").Value = 0.04544590891535629
$ws.Range("T15").Value = 0.04544590891535629
$ws.Range("G16").Value = 12.45667266666667
$ws.Range("H16").Value = 37.370018
$ws.Range("I16").Value = 0.2744805580967825
$ws.Range("J16").Value = 0.2744805580967826
$ws.Range("M16").Value = 7.682722666666667
$ws.Range("N16").Value = 23.048168
$ws.Range("O16").Value = 0.5166467859527435
$ws.Range("P16").Value = 0.5166467859527435
$ws.Range("Q16").Value = 95.70116144744712
$ws.Range("R16").Value = 861.3104530270241
$ws.Range("S16").Value = 0.141809498147218
$ws.Range("T16").Value = 0.141809498147218
$ws.Range("G17").Value = 12.45667266666667
$ws.Range("H17").Value = 37.370018
$ws.Range("I17").Value = 0.2744805580967825
$ws.Range("J17").Value = 0.2744805580967826
$ws.Range("M17").Value = 3.105673
$ws.Range("N17").Value = 9.317019
$ws.Range("O17").Value = 0.2088499147095181
$ws.Range("P17").Value = 0.2088499147095181
$ws.Range("Q17").Value = 38.68635197070467
$ws.Range("R17").Value = 348.177167736342
$ws.Range("S17").Value = 0.05732524114793395
$ws.Range("T17").Value = 0.05732524114793396
